# Update TPM-derived values in sheet1 per new data (commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 99.28451533333333
$ws.Range("H2").Value = 297.853546
$ws.Range("I2").Value = 0.02270354261926982
$ws.Range("J2").Value = 0.02270354261926982
$ws.Range("M2").Value = 6.045145666666667
$ws.Range("N2").Value = 18.135437
$ws.Range("O2").Value = 0.8160840232643366
$ws.Range("P2").Value = 0.8160840232643367
$ws.Range("Q2").Value = 600.1893576344003
$ws.Range("R2").Value = 5401.704218709603
$ws.Range("S2").Value = 0.01852799840308705
$ws.Range("T2").Value = 0.01852799840308705
$ws.Range("G3").Value = 99.28451533333333
$ws.Range("H3").Value = 297.853546
$ws.Range("I3").Value = 0.02270354261926982
$ws.Range("J3").Value = 0.02270354261926982
$ws.Range("O3").Value = 0.09212864864242169
$ws.Range("P3").Value = 0.09212864864242169
$ws.Range("Q3").Value = 67.75605559246445
$ws.Range("R3").Value = 609.8045003321801
$ws.Range("S3").Value = 0.002091646700908956
$ws.Range("T3").Value = 0.002091646700908955
$ws.Range("G4").Value = 99.28451533333333
$ws.Range("H4").Value = 297.853546
$ws.Range("I4").Value = 0.02270354261926982
$ws.Range("J4").Value = 0.02270354261926982
$ws.Range("M4").Value = 0.6799149999999999
$ws.Range("N4").Value = 2.039745
$ws.Range("O4").Value = 0.09178732809324164
$ws.Range("P4").Value = 0.09178732809324165
$ws.Range("Q4").Value = 67.50503124286332
$ws.Range("R4").Value = 607.5452811857699
$ws.Range("S4").Value = 0.002083897515273814
$ws.Range("T4").Value = 0.002083897515273814
$ws.Range("I5").Value = 0.9171714767027319
$ws.Range("J5").Value = 0.9171714767027318
$ws.Range("M5").Value = 6.045145666666667
$ws.Range("N5").Value = 18.135437
$ws.Range("O5").Value = 0.8160840232643366
$ws.Range("P5").Value = 0.8160840232643367
$ws.Range("Q5").Value = 24246.28476155018
$ws.Range("R5").Value = 218216.5628539516
$ws.Range("S5").Value = 0.7484889887308582
$ws.Range("T5").Value = 0.7484889887308582
$ws.Range("I6").Value = 0.9171714767027319
$ws.Range("J6").Value = 0.9171714767027318
$ws.Range("O6").Value = 0.09212864864242169
$ws.Range("P6").Value = 0.09212864864242169
$ws.Range("S6").Value = 0.08449776872199705
$ws.Range("T6").Value = 0.08449776872199703
$ws.Range("I7").Value = 0.9171714767027319
$ws.Range("J7").Value = 0.9171714767027318
$ws.Range("M7").Value = 0.6799149999999999
$ws.Range("N7").Value = 2.039745
$ws.Range("O7").Value = 0.09178732809324164
$ws.Range("P7").Value = 0.09178732809324165
$ws.Range("Q7").Value = 2727.049704451465
$ws.Range("R7").Value = 24543.44734006319
$ws.Range("S7").Value = 0.08418471924987658
$ws.Range("T7").Value = 0.08418471924987658
$ws.Range("G8").Value = 249.2612966666667
$ws.Range("H8").Value = 747.7838899999999
$ws.Range("I8").Value = 0.05699896356653876
$ws.Range("J8").Value = 0.05699896356653875
$ws.Range("M8").Value = 6.045145666666667
$ws.Range("N8").Value = 18.135437
$ws.Range("O8").Value = 0.8160840232643366
$ws.Range("P8").Value = 0.8160840232643367
$ws.Range("Q8").Value = 1506.820847412215
$ws.Range("R8").Value = 13561.38762670993
$ws.Range("S8").Value = 0.04651594350927829
$ws.Range("T8").Value = 0.04651594350927829
$ws.Range("G9").Value = 249.2612966666667
$ws.Range("H9").Value = 747.7838899999999
$ws.Range("I9").Value = 0.05699896356653876
$ws.Range("J9").Value = 0.05699896356653875
$ws.Range("O9").Value = 0.09212864864242169
$ws.Range("P9").Value = 0.09212864864242169
$ws.Range("Q9").Value = 170.1067101681889
$ws.Range("R9").Value = 1530.9603915137
$ws.Range("S9").Value = 0.005251237487403844
$ws.Range("T9").Value = 0.005251237487403844
$ws.Range("G10").Value = 249.2612966666667
$ws.Range("H10").Value = 747.7838899999999
$ws.Range("I10").Value = 0.05699896356653876
$ws.Range("J10").Value = 0.05699896356653875
$ws.Range("M10").Value = 0.6799149999999999
$ws.Range("N10").Value = 2.039745
$ws.Range("O10").Value = 0.09178732809324164
$ws.Range("P10").Value = 0.09178732809324165
$ws.Range("Q10").Value = 169.4764945231166
$ws.Range("R10").Value = 1525.28845070805
$ws.Range("S10").Value = 0.00523178256985662
$ws.Range("T10").Value = 0.00523178256985662
$ws.Range("G11").Value = 13.67033766666667
$ws.Range("H11").Value = 41.011013
$ws.Range("I11").Value = 0.003126017111459632
$ws.Range("J11").Value = 0.003126017111459632
$ws.Range("M11").Value = 6.045145666666667
$ws.Range("N11").Value = 18.135437
$ws.Range("O11").Value = 0.8160840232643366
$ws.Range("P11").Value = 0.8160840232643367
$ws.Range("Q11").Value = 82.63918250752012
$ws.Range("R11").Value = 743.7526425676811
$ws.Range("S11").Value = 0.002551092621113137
$ws.Range("T11").Value = 0.002551092621113137
$ws.Range("G12").Value = 13.67033766666667
$ws.Range("H12").Value = 41.011013
$ws.Range("I12").Value = 0.003126017111459632
$ws.Range("J12").Value = 0.003126017111459632
$ws.Range("O12").Value = 0.09212864864242169
$ws.Range("P12").Value = 0.09212864864242169
$ws.Range("Q12").Value = 9.329230805032223
$ws.Range("R12").Value = 83.96307724528999
$ws.Range("S12").Value = 0.0002879957321118624
$ws.Range("T12").Value = 0.0002879957321118624
$ws.Range("G13").Value = 13.67033766666667
$ws.Range("H13").Value = 41.011013
$ws.Range("I13").Value = 0.003126017111459632
$ws.Range("J13").Value = 0.003126017111459632
$ws.Range("M13").Value = 0.6799149999999999
$ws.Range("N13").Value = 2.039745
$ws.Range("O13").Value = 0.09178732809324164
$ws.Range("P13").Value = 0.09178732809324165
$ws.Range("Q13").Value = 9.294667634631665
$ws.Range("R13").Value = 83.652008711685
$ws.Range("S13").Value = 0.0002869287582346328
$ws.Range("T13").Value = 0.0002869287582346328

Write-Output "Updated cells with new TPM values"
